# Mark the "Invalid Login" (row 7) and "Empty Fields" (row 8) test cases as
# completed by Joel, matching the existing "Completion" entries already set
# for the "Valid Login" test cases in rows 5 and 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Invalid Login
$ws.Range("D7").Value = "Joel"
$ws.Range("D5").Copy()
$ws.Range("D7").PasteSpecial(-4122)

# Row 8 - Empty Fields
$ws.Range("D8").Value = "Joel"
$ws.Range("D8").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Reflect the reviewer's scroll position / selection at the time of the edit.
$win = $excel.ActiveWindow
$null = ($win.ScrollRow = 4)
$null = ($win.ScrollColumn = 1)
$null = $ws.Range("G10").Select()
